# data cleanup continued in multi dfs, major code refactoring taking place contiued
#
# The sheet is a sorted (by player name) summary of "player -> award count".
# Two new players need to be inserted into their correct alphabetically-sorted
# rows ("Leandro Barbosa" and "Malcolm Brogdon"), which shifts every
# subsequent row down. We rebuild the full, final sorted A2:B36 range in one
# pass rather than trying to patch individual shifted rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    'Aaron McKie',
    'Antawn Jamison',
    'Anthony Mason',
    'Ben Gordon',
    'Bill Walton',
    'Bobby Jackson',
    'Bobby Jones',
    'Clifford Robinson',
    'Corliss Williamson',
    'Danny Manning',
    'Darrell Armstrong',
    'Dell Curry',
    'Detlef Schrempf',
    'Eddie Johnson',
    'Eric Gordon',
    'J.R. Smith',
    'Jamal Crawford',
    'James Harden',
    'Jason Terry',
    'John Starks',
    'Jordan Clarkson',
    'Kevin McHale',
    'Lamar Odom',
    'Leandro Barbosa',
    'Lou Williams',
    'Malcolm Brogdon',
    'Manu Ginóbili',
    'Mike Miller',
    'Montrezl Harrell',
    'Naz Reid',
    'Ricky Pierce',
    'Rodney Rogers',
    'Roy Tarpley',
    'Toni Kukoč',
    'Tyler Herro'
)

$counts = @(1,1,1,1,1,1,1,1,1,1,1,1,2,1,1,1,3,1,1,1,1,2,1,1,3,1,1,1,1,1,2,1,1,1,1)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value2 = $names[$i]
    $ws.Cells.Item($row, 2).Value2 = $counts[$i]
}
